# Add a new "2022-Q3" sheet (right after "总计") with its fund-holding
# data, and add the corresponding summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by copying the existing
#    "2022-Q2" sheet (so it inherits the same sheetPr / styles / page
#    setup), placed right after "总计".
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$srcSheet.Copy($null, $totalSheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# Give rows 3 and 4 the same "index column" style (A2's style, s=2)
# that row 2 already has, before we fill in the extra fund rows.
$q3Sheet.Cells.Item(2, 1).Copy()
$q3Sheet.Cells.Item(3, 1).PasteSpecial(-4122)
$q3Sheet.Cells.Item(4, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2: 001628 招商体育文化休闲股票A
$q3Sheet.Cells.Item(2, 1).Value = 0
$q3Sheet.Cells.Item(2, 2).Style = "Normal"
$q3Sheet.Cells.Item(2, 2).Value = "'001628"
$q3Sheet.Cells.Item(2, 2).Style = "Normal"
$q3Sheet.Cells.Item(2, 3).Value = "招商体育文化休闲股票A"
$q3Sheet.Cells.Item(2, 4).Style = "Normal"
$q3Sheet.Cells.Item(2, 4).Value = "'2.23"
$q3Sheet.Cells.Item(2, 4).Style = "Normal"
$q3Sheet.Cells.Item(2, 5).Style = "Normal"
$q3Sheet.Cells.Item(2, 5).Value = "'92.42"
$q3Sheet.Cells.Item(2, 5).Style = "Normal"
$q3Sheet.Cells.Item(2, 6).Style = "Normal"
$q3Sheet.Cells.Item(2, 6).Value = "'5.46"
$q3Sheet.Cells.Item(2, 6).Style = "Normal"
$q3Sheet.Cells.Item(2, 7).Style = "Normal"
$q3Sheet.Cells.Item(2, 7).Value = "'0.1218"
$q3Sheet.Cells.Item(2, 7).Style = "Normal"
$q3Sheet.Cells.Item(2, 8).Value = 1

# Row 3: 003397 银华体育文化灵活配置混合
$q3Sheet.Cells.Item(3, 1).Value = 1
$q3Sheet.Cells.Item(3, 2).Style = "Normal"
$q3Sheet.Cells.Item(3, 2).Value = "'003397"
$q3Sheet.Cells.Item(3, 2).Style = "Normal"
$q3Sheet.Cells.Item(3, 3).Value = "银华体育文化灵活配置混合"
$q3Sheet.Cells.Item(3, 4).Style = "Normal"
$q3Sheet.Cells.Item(3, 4).Value = "'0.32"
$q3Sheet.Cells.Item(3, 4).Style = "Normal"
$q3Sheet.Cells.Item(3, 5).Style = "Normal"
$q3Sheet.Cells.Item(3, 5).Value = "'81.07"
$q3Sheet.Cells.Item(3, 5).Style = "Normal"
$q3Sheet.Cells.Item(3, 6).Style = "Normal"
$q3Sheet.Cells.Item(3, 6).Value = "'4.46"
$q3Sheet.Cells.Item(3, 6).Style = "Normal"
$q3Sheet.Cells.Item(3, 7).Style = "Normal"
$q3Sheet.Cells.Item(3, 7).Value = "'0.0143"
$q3Sheet.Cells.Item(3, 7).Style = "Normal"
$q3Sheet.Cells.Item(3, 8).Value = 7

# Row 4: 015395 招商体育文化休闲股票C
$q3Sheet.Cells.Item(4, 1).Value = 2
$q3Sheet.Cells.Item(4, 2).Style = "Normal"
$q3Sheet.Cells.Item(4, 2).Value = "'015395"
$q3Sheet.Cells.Item(4, 2).Style = "Normal"
$q3Sheet.Cells.Item(4, 3).Value = "招商体育文化休闲股票C"
$q3Sheet.Cells.Item(4, 4).Style = "Normal"
$q3Sheet.Cells.Item(4, 4).Value = "'0.25"
$q3Sheet.Cells.Item(4, 4).Style = "Normal"
$q3Sheet.Cells.Item(4, 5).Style = "Normal"
$q3Sheet.Cells.Item(4, 5).Value = "'92.42"
$q3Sheet.Cells.Item(4, 5).Style = "Normal"
$q3Sheet.Cells.Item(4, 6).Style = "Normal"
$q3Sheet.Cells.Item(4, 6).Value = "'5.46"
$q3Sheet.Cells.Item(4, 6).Style = "Normal"
$q3Sheet.Cells.Item(4, 7).Style = "Normal"
$q3Sheet.Cells.Item(4, 7).Value = "'0.0136"
$q3Sheet.Cells.Item(4, 7).Style = "Normal"
$q3Sheet.Cells.Item(4, 8).Value = 1

# ---------------------------------------------------------------------
# 2) Update the "总计" (total) sheet: prepend a new "2022-Q3" summary
#    row and push the previous rows down by one, re-numbering the
#    index column (A) sequentially.
# ---------------------------------------------------------------------

# Give the new row 6 the same "index column" style (s=2) as the rows
# above it before writing into it.
$totalSheet.Cells.Item(5, 1).Copy()
$totalSheet.Cells.Item(6, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 3
$totalSheet.Cells.Item(2, 4).Value = 0.15

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 1
$totalSheet.Cells.Item(3, 4).Value = 0.02

$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 4
$totalSheet.Cells.Item(4, 4).Value = 0.58

$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q4"
$totalSheet.Cells.Item(5, 3).Value = 5
$totalSheet.Cells.Item(5, 4).Value = 1.49

$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(6, 2).Value = "2021-Q2"
$totalSheet.Cells.Item(6, 3).Value = 1
$totalSheet.Cells.Item(6, 4).Value = 0.02

# ---------------------------------------------------------------------
# 3) Restore the originally-active tab ("2021-Q2", now the last sheet).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
